$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: add empty styled cell R3 (matches Q3's style) ---
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 4: add R4 = 2021 (matches Q4's style) ---
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R4").Value = 2021

# --- Row 5: R5 gets a value and a new "vertical top" direct format ---
$ws.Range("R5").Value = 0.9
$ws.Range("R5").VerticalAlignment = -4160   # xlVAlignTop

# --- Row 6: R6 = 6.5 (matches Q6's style) ---
$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R6").Value = 6.5

# --- Update the active selection shown in the sheet view ---
$null = $ws.Range("T5").Select()

Write-Host "done"
